$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 35
$ws.Range("J6").Value = 63
$ws.Range("L6").Value = 189
$ws.Range("N6").Value = -413
$ws.Range("H11").Value = 430.6154
$ws.Range("I11").Value = 430.6154
$ws.Range("K11").Value = 430.6154
$ws.Range("M11").Value = -290.6154
$ws.Range("H17").Value = 3099.9
$ws.Range("J17").Value = 3333
$ws.Range("L17").Value = 9999
$ws.Range("N17").Value = -10335
$ws.Range("H38").Value = 108
$ws.Range("I38").Value = 18.545454
$ws.Range("J38").Value = 600
$ws.Range("K38").Value = 55.636362
$ws.Range("L38").Value = 1800
$ws.Range("M38").Value = 316.363638
$ws.Range("N38").Value = -2544
$ws.Range("H41").Value = 542.75
$ws.Range("I41").Value = 542.75
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 542.75
$ws.Range("L41").Value = 0
$ws.Range("M41").Value = -102.75
$ws.Range("N41").ClearContents()
$ws.Range("H58").Value = 446.6
$ws.Range("I58").Value = 11
$ws.Range("J58").Value = 1100
$ws.Range("K58").Value = 33
$ws.Range("L58").Value = 3300
$ws.Range("M58").Value = 117
$ws.Range("N58").Value = -3600
$ws.Range("H103").Value = 1680.6666
$ws.Range("I103").Value = 1450
$ws.Range("K103").Value = 4350
$ws.Range("M103").Value = -3764
$ws.Range("H109").Value = 0
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("M109").ClearContents()
$ws.Range("N109").ClearContents()
$ws.Range("H112").Value = 1750
$ws.Range("I112").Value = 1500
$ws.Range("J112").Value = 2000
$ws.Range("K112").Value = 4500
$ws.Range("L112").Value = 6000
$ws.Range("M112").Value = -3392
$ws.Range("N112").Value = -8216
$ws.Range("H125").Value = 980
$ws.Range("I125").Value = 980
$ws.Range("K125").Value = 8820
$ws.Range("M125").Value = -6360
$ws.Range("H131").Value = 800.6667
$ws.Range("I131").Value = 800.6667
$ws.Range("K131").Value = 2402.0001
$ws.Range("M131").Value = 2637.9999
$ws.Range("H132").Value = 112055.5
$ws.Range("I132").Value = 133661.8
$ws.Range("K132").Value = 400985.4
$ws.Range("M132").Value = -398455.4

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("M61").ClearContents()
$ws.Range("H102").Value = 26376374
$ws.Range("I102").Value = 1572999.8
$ws.Range("K102").Value = 1572999.8
$ws.Range("M102").Value = -1571377.8
$ws.Range("H122").Value = 1674.75
$ws.Range("I122").Value = 1449.5
$ws.Range("K122").Value = 4348.5
$ws.Range("M122").Value = -1898.5
$ws.Range("H132").Value = 3974.6667
$ws.Range("I132").Value = 3974.6667
$ws.Range("K132").Value = 11924.0001
$ws.Range("M132").Value = -9394.000100000001
$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").ClearContents()

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 520.8333
$ws.Range("I7").Value = 453
$ws.Range("J7").Value = 588.6667
$ws.Range("K7").Value = 453
$ws.Range("L7").Value = 588.6667
$ws.Range("M7").Value = -340
$ws.Range("N7").Value = -814.6667
$ws.Range("H80").Value = 1398.3334
$ws.Range("I80").Value = 1397.5
$ws.Range("J80").Value = 1400
$ws.Range("K80").Value = 1397.5
$ws.Range("L80").Value = 1400
$ws.Range("M80").Value = -399.5
$ws.Range("N80").Value = -3396
$ws.Range("H83").Value = 1398.3334
$ws.Range("I83").Value = 1397.5
$ws.Range("J83").Value = 1400
$ws.Range("K83").Value = 6987.5
$ws.Range("L83").Value = 7000
$ws.Range("M83").Value = -1995.5
$ws.Range("N83").Value = -16984

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 3794.375
$ws.Range("I19").Value = 899.75
$ws.Range("J19").Value = 6689
$ws.Range("K19").Value = 899.75
$ws.Range("L19").Value = 6689
$ws.Range("M19").Value = -729.75
$ws.Range("N19").Value = -7029
$ws.Range("H21").Value = 5000
$ws.Range("J21").Value = 5000
$ws.Range("L21").Value = 5000
$ws.Range("N21").Value = -5470
$ws.Range("H24").Value = 3794.375
$ws.Range("I24").Value = 899.75
$ws.Range("J24").Value = 6689
$ws.Range("K24").Value = 899.75
$ws.Range("L24").Value = 6689
$ws.Range("M24").Value = -729.75
$ws.Range("N24").Value = -7029
$ws.Range("H25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("N25").ClearContents()
$ws.Range("H31").Value = 6771.4
$ws.Range("I31").Value = 3987.5
$ws.Range("J31").Value = 8627.333000000001
$ws.Range("K31").Value = 3987.5
$ws.Range("L31").Value = 8627.333000000001
$ws.Range("M31").Value = -3692.5
$ws.Range("N31").Value = -9217.333000000001
$ws.Range("H34").Value = 6771.4
$ws.Range("I34").Value = 3987.5
$ws.Range("J34").Value = 8627.333000000001
$ws.Range("K34").Value = 3987.5
$ws.Range("L34").Value = 8627.333000000001
$ws.Range("M34").Value = -3785.5
$ws.Range("N34").Value = -9031.333000000001
$ws.Range("H43").Value = 22166.666
$ws.Range("J43").Value = 22166.666
$ws.Range("L43").Value = 22166.666
$ws.Range("N43").Value = -22534.666
$ws.Range("H101").Value = 22166.666
$ws.Range("J101").Value = 22166.666
$ws.Range("L101").Value = 22166.666
$ws.Range("N101").Value = -28656.666

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("M8").ClearContents()
$ws.Range("H68").Value = 1815.7778
$ws.Range("I68").Value = 2192.25
$ws.Range("J68").Value = 1514.6
$ws.Range("K68").Value = 6576.75
$ws.Range("L68").Value = 4543.799999999999
$ws.Range("M68").Value = -5765.75
$ws.Range("N68").Value = -6165.799999999999
$ws.Range("H71").Value = 1815.7778
$ws.Range("I71").Value = 2192.25
$ws.Range("J71").Value = 1514.6
$ws.Range("K71").Value = 19730.25
$ws.Range("L71").Value = 13631.4
$ws.Range("M71").Value = -15674.25
$ws.Range("N71").Value = -21743.4

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 83336750
$ws.Range("I70").Value = 3998.3333
$ws.Range("J70").Value = 166669500
$ws.Range("K70").Value = 3998.3333
$ws.Range("L70").Value = 166669500
$ws.Range("M70").Value = -3728.3333
$ws.Range("N70").Value = -166670040
$ws.Range("H73").Value = 83336750
$ws.Range("I73").Value = 3998.3333
$ws.Range("J73").Value = 166669500
$ws.Range("K73").Value = 3998.3333
$ws.Range("L73").Value = 166669500
$ws.Range("M73").Value = -3062.3333
$ws.Range("N73").Value = -166671372
$ws.Range("H122").Value = 4671.35
$ws.Range("I122").Value = 4301.5
$ws.Range("K122").Value = 12904.5
$ws.Range("M122").Value = -10454.5

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1797.75
$ws.Range("I61").Value = 1797.75
$ws.Range("K61").Value = 1797.75
$ws.Range("M61").Value = -1595.75
$ws.Range("H82").Value = 3178
$ws.Range("I82").Value = 1773.375
$ws.Range("K82").Value = 1773.375
$ws.Range("M82").Value = -1412.375
$ws.Range("H85").Value = 3178
$ws.Range("I85").Value = 1773.375
$ws.Range("K85").Value = 1773.375
$ws.Range("M85").Value = -525.375
$ws.Range("H103").Value = 14200
$ws.Range("J103").Value = 14200
$ws.Range("L103").Value = 14200
$ws.Range("N103").Value = -16544
$ws.Range("H113").Value = 1797.75
$ws.Range("I113").Value = 1797.75
$ws.Range("K113").Value = 1797.75
$ws.Range("M113").Value = 372.25
$ws.Range("H122").Value = 3696.4075
$ws.Range("I122").Value = 3141.9167
$ws.Range("K122").Value = 9425.750100000001
$ws.Range("M122").Value = -6975.750100000001
$ws.Range("H132").Value = 2878.75
$ws.Range("J132").Value = 2938.3333
$ws.Range("L132").Value = 8814.999899999999
$ws.Range("N132").Value = -13874.9999

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 284.85715
$ws.Range("I107").Value = 299.33334
$ws.Range("K107").Value = 898.0000200000001
$ws.Range("M107").Value = 1021.99998
$ws.Range("H126").Value = 5175
$ws.Range("I126").Value = 4432.6665
$ws.Range("K126").Value = 13297.9995
$ws.Range("M126").Value = -10827.9995
